# Auto-generated edit script: updates crypto Price (D) and Volume(1h) (E) columns
# for rows 2-51 per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each touched cell to Text format before assigning, so that purely
# numeric-looking strings (e.g. "246.71", "0.0951") are preserved exactly as
# text instead of being coerced into floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.891.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.234.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.79"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.80%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.53"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.33%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.566.73"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.81"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.240.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.793.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.78%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.11"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.98"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.08"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +10.85%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.70"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.51"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.78"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0830"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.121"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.53"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0300"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.46"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.89"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.18"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.68"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.92%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.73"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.68"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.76%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.996"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -12.25%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.17"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.19%  "
